# Applies row-level financial data updates across the Leviathan Profits workbook
# (one table per crafting-class worksheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 116778
$ws.Range("I62").Value = 253800.75
$ws.Range("J62").Value = 7159.8
$ws.Range("K62").Value = 253800.75
$ws.Range("L62").Value = 7159.8
$ws.Range("M62").Value = -253176.75
$ws.Range("N62").Value = -8407.799999999999
# Row 65
$ws.Range("H65").Value = 116778
$ws.Range("I65").Value = 253800.75
$ws.Range("J65").Value = 7159.8
$ws.Range("K65").Value = 1269003.75
$ws.Range("L65").Value = 35799
$ws.Range("M65").Value = -1265883.75
$ws.Range("N65").Value = -42039
# Row 106
$ws.Range("H106").Value = 10520.692
$ws.Range("I106").Value = 2199
$ws.Range("J106").Value = 15721.75
$ws.Range("K106").Value = 2199
$ws.Range("L106").Value = 15721.75
$ws.Range("M106").Value = -1568
$ws.Range("N106").Value = -16983.75
# Row 134
$ws.Range("H134").Value = 119824
$ws.Range("J134").Value = 99932.664
$ws.Range("L134").Value = 99932.664
$ws.Range("N134").Value = -110072.664

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 36490.05
$ws.Range("J32").Value = 210979.6
$ws.Range("L32").Value = 210979.6
$ws.Range("N32").Value = -211553.6
# Row 45
$ws.Range("H45").Value = 845245.75
$ws.Range("I45").Value = 1445474.1
$ws.Range("K45").Value = 1445474.1
$ws.Range("M45").Value = -1445097.1
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 89999
$ws.Range("J135").Value = 89999
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139
# Row 137
$ws.Range("H137").Value = 88886.664
$ws.Range("J137").Value = 88886.664
$ws.Range("L137").Value = 88886.664
$ws.Range("N137").Value = -99086.664
# Row 138
$ws.Range("H138").Value = 111249.25
$ws.Range("J138").Value = 111249.25
$ws.Range("L138").Value = 111249.25
$ws.Range("N138").Value = -121529.25
# Row 140
$ws.Range("H140").Value = 139994
$ws.Range("J140").Value = 139994
$ws.Range("L140").Value = 139994
$ws.Range("N140").Value = -150354

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 25996
$ws.Range("J81").Value = 25996
$ws.Range("L81").Value = 25996
$ws.Range("N81").Value = -28118
# Row 84
$ws.Range("H84").Value = 25996
$ws.Range("J84").Value = 25996
$ws.Range("L84").Value = 77988
$ws.Range("N84").Value = -88596
# Row 105
$ws.Range("H105").Value = 2634754.8
$ws.Range("I105").Value = 2860333.8
$ws.Range("K105").Value = 2860333.8
$ws.Range("M105").Value = -2858586.8
# Row 134
$ws.Range("H134").Value = 1582.3572
$ws.Range("I134").Value = 1582.3572
$ws.Range("K134").Value = 4747.071599999999
$ws.Range("M134").Value = -2212.071599999999
# Row 139
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280
# Row 140
$ws.Range("H140").Value = 53568.75
$ws.Range("J140").Value = 53568.75
$ws.Range("L140").Value = 53568.75
$ws.Range("N140").Value = -63928.75
# Row 141
$ws.Range("H141").Value = 92138
$ws.Range("J141").Value = 99995.25
$ws.Range("L141").Value = 99995.25
$ws.Range("N141").Value = -110355.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2337.0732
$ws.Range("J31").Value = 3165.9375
$ws.Range("L31").Value = 3165.9375
$ws.Range("N31").Value = -3755.9375
# Row 34
$ws.Range("H34").Value = 2337.0732
$ws.Range("J34").Value = 3165.9375
$ws.Range("L34").Value = 3165.9375
$ws.Range("N34").Value = -3569.9375
# Row 133
$ws.Range("H133").Value = 49997.5
$ws.Range("J133").Value = 49997.5
$ws.Range("L133").Value = 49997.5
$ws.Range("N133").Value = -55057.5
# Row 134
$ws.Range("H134").Value = 3416.7856
$ws.Range("I134").Value = 3281.4092
$ws.Range("J134").Value = 3913.1667
$ws.Range("K134").Value = 9844.2276
$ws.Range("L134").Value = 11739.5001
$ws.Range("M134").Value = -7309.2276
$ws.Range("N134").Value = -16809.5001
# Row 141
$ws.Range("H141").Value = 328499.66
$ws.Range("J141").Value = 328499.66
$ws.Range("L141").Value = 328499.66
$ws.Range("N141").Value = -338859.66

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 177.3158
$ws.Range("J17").Value = 235.44444
$ws.Range("L17").Value = 706.33332
$ws.Range("N17").Value = -1044.33332
# Row 34
$ws.Range("H34").Value = 1665.4286
$ws.Range("J34").Value = 3599.3333
$ws.Range("L34").Value = 10797.9999
$ws.Range("N34").Value = -10965.9999
# Row 38
$ws.Range("H38").Value = 156
$ws.Range("I38").Value = 93.333336
$ws.Range("K38").Value = 280.000008
$ws.Range("M38").Value = 66.99999200000002
# Row 39
$ws.Range("H39").Value = 4842.857
$ws.Range("J39").Value = 9966.666999999999
$ws.Range("L39").Value = 29900.001
$ws.Range("N39").Value = -30488.001
# Row 55
$ws.Range("H55").Value = 62504950
$ws.Range("J55").Value = 62504950
$ws.Range("L55").Value = 187514850
$ws.Range("N55").Value = -187515204
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 128
$ws.Range("H128").Value = 246247.86
$ws.Range("I128").Value = 246247.86
$ws.Range("K128").Value = 738743.58
$ws.Range("M128").Value = -733763.58

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3502.2222
$ws.Range("I126").Value = 3620.1667
$ws.Range("K126").Value = 10860.5001
$ws.Range("M126").Value = -8390.500100000001
# Row 135
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140
# Row 140
$ws.Range("H140").Value = 76497.664
$ws.Range("J140").Value = 76497.664
$ws.Range("L140").Value = 76497.664
$ws.Range("N140").Value = -86857.664

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1849.6666
$ws.Range("I22").Value = 1819.6
$ws.Range("K22").Value = 1819.6
$ws.Range("M22").Value = -1524.6
# Row 27
$ws.Range("H27").Value = 1849.6666
$ws.Range("I27").Value = 1819.6
$ws.Range("K27").Value = 1819.6
$ws.Range("M27").Value = -1712.6
# Row 45
$ws.Range("H45").Value = 20633
$ws.Range("I45").Value = 15012.333
$ws.Range("K45").Value = 15012.333
$ws.Range("M45").Value = -14605.333
# Row 133
$ws.Range("H133").Value = 78999
$ws.Range("J133").Value = 78999
$ws.Range("L133").Value = 78999
$ws.Range("N133").Value = -84059
# Row 138
$ws.Range("H138").Value = 73000
$ws.Range("J138").Value = 73000
$ws.Range("L138").Value = 73000
$ws.Range("N138").Value = -83280

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1914.2858
$ws.Range("I96").Value = 1880
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1880
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -507
$ws.Range("N96").Value = -4746

Write-Output "Applied 40 row updates across 8 sheets."
